# Auto-generated Excel COM-interop script
# Updates FFXIV Excalibur market-board profit data (currentAveragePrice /
# LevePriceNQ/HQ / LeveProfitNQ/HQ columns) on each job sheet, per the
# scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H39").Value = 604.2857
$ws.Range("I39").Value = 494.81818
$ws.Range("K39").Value = 1484.45454
$ws.Range("M39").Value = -1188.45454
$ws.Range("H40").Value = 3359.5386
$ws.Range("I40").Value = 3109.5
$ws.Range("K40").Value = 3109.5
$ws.Range("M40").Value = -2934.5
$ws.Range("H82").Value = 1078.1428
$ws.Range("I82").Value = 1078.1428
$ws.Range("K82").Value = 3234.4284
$ws.Range("M82").Value = -2828.4284
$ws.Range("H85").Value = 1078.1428
$ws.Range("I85").Value = 1078.1428
$ws.Range("K85").Value = 3234.4284
$ws.Range("M85").Value = -1830.4284
$ws.Range("H86").Value = 13592150
$ws.Range("I86").Value = 7819178
$ws.Range("K86").Value = 7819178
$ws.Range("M86").Value = -7818055
$ws.Range("H89").Value = 13592150
$ws.Range("I89").Value = 7819178
$ws.Range("K89").Value = 39095890
$ws.Range("M89").Value = -39090274
$ws.Range("H111").Value = 2976.875
$ws.Range("I111").Value = 2173.2
$ws.Range("J111").Value = 4316.3335
$ws.Range("K111").Value = 6519.599999999999
$ws.Range("L111").Value = 12949.0005
$ws.Range("M111").Value = -3452.599999999999
$ws.Range("N111").Value = -19083.0005
$ws.Range("H137").Value = 2445.0625
$ws.Range("I137").Value = 977.6875
$ws.Range("J137").Value = 3912.4375
$ws.Range("K137").Value = 2933.0625
$ws.Range("L137").Value = 11737.3125
$ws.Range("M137").Value = -383.0625
$ws.Range("N137").Value = -16837.3125
$ws.Range("H138").Value = 4422.074
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4422.074
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 13266.222
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -23546.222

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H5").Value = 237.83333
$ws.Range("I5").Value = 244.5
$ws.Range("J5").Value = 224.5
$ws.Range("K5").Value = 244.5
$ws.Range("L5").Value = 224.5
$ws.Range("M5").Value = -132.5
$ws.Range("N5").Value = -448.5
$ws.Range("H32").Value = 796.2394399999999
$ws.Range("I32").Value = 721.9
$ws.Range("K32").Value = 721.9
$ws.Range("M32").Value = -434.9
$ws.Range("H63").Value = 5580.6875
$ws.Range("J63").Value = 6499.4614
$ws.Range("L63").Value = 6499.4614
$ws.Range("N63").Value = -7871.4614
$ws.Range("H66").Value = 5580.6875
$ws.Range("J66").Value = 6499.4614
$ws.Range("L66").Value = 32497.307
$ws.Range("N66").Value = -39361.307
$ws.Range("H88").Value = 1518.3077
$ws.Range("I88").Value = 1313.25
$ws.Range("J88").Value = 1609.4445
$ws.Range("K88").Value = 1313.25
$ws.Range("L88").Value = 1609.4445
$ws.Range("M88").Value = -907.25
$ws.Range("N88").Value = -2421.4445
$ws.Range("H91").Value = 1518.3077
$ws.Range("I91").Value = 1313.25
$ws.Range("J91").Value = 1609.4445
$ws.Range("K91").Value = 1313.25
$ws.Range("L91").Value = 1609.4445
$ws.Range("M91").Value = 90.75
$ws.Range("N91").Value = -4417.4445
$ws.Range("H110").Value = 2033.3914
$ws.Range("I110").Value = 1248.55
$ws.Range("K110").Value = 1248.55
$ws.Range("M110").Value = 796.45
$ws.Range("H122").Value = 2822.4062
$ws.Range("I122").Value = 2262.1177
$ws.Range("J122").Value = 3457.4
$ws.Range("K122").Value = 6786.353099999999
$ws.Range("L122").Value = 10372.2
$ws.Range("M122").Value = -4336.353099999999
$ws.Range("N122").Value = -15272.2

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H4").Value = 237.83333
$ws.Range("I4").Value = 244.5
$ws.Range("J4").Value = 224.5
$ws.Range("K4").Value = 244.5
$ws.Range("L4").Value = 224.5
$ws.Range("M4").Value = -129.5
$ws.Range("N4").Value = -454.5
$ws.Range("H105").Value = 3000.8096
$ws.Range("I105").Value = 3173.7222
$ws.Range("K105").Value = 3173.7222
$ws.Range("M105").Value = -1426.7222
$ws.Range("H107").Value = 1550.6383
$ws.Range("I107").Value = 1935.8846
$ws.Range("J107").Value = 1073.6666
$ws.Range("K107").Value = 1935.8846
$ws.Range("L107").Value = 1073.6666
$ws.Range("M107").Value = -15.88460000000009
$ws.Range("N107").Value = -4913.6666
$ws.Range("H134").Value = 1794.45
$ws.Range("I134").Value = 836.2632
$ws.Range("K134").Value = 2508.7896
$ws.Range("M134").Value = 26.21039999999994

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 16913.45
$ws.Range("J31").Value = 38959.332
$ws.Range("L31").Value = 38959.332
$ws.Range("N31").Value = -39549.332
$ws.Range("H34").Value = 16913.45
$ws.Range("J34").Value = 38959.332
$ws.Range("L34").Value = 38959.332
$ws.Range("N34").Value = -39363.332
$ws.Range("H122").Value = 1978
$ws.Range("I122").Value = 1984.6
$ws.Range("J122").Value = 1953.25
$ws.Range("K122").Value = 5953.799999999999
$ws.Range("L122").Value = 5859.75
$ws.Range("M122").Value = -3503.799999999999
$ws.Range("N122").Value = -10759.75

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H12").Value = 273.48
$ws.Range("I12").Value = 36.666668
$ws.Range("K12").Value = 110.000004
$ws.Range("M12").Value = 62.999996
$ws.Range("H17").Value = 200000200
$ws.Range("I17").Value = 32.5
$ws.Range("J17").Value = 333333630
$ws.Range("K17").Value = 97.5
$ws.Range("L17").Value = 1000000890
$ws.Range("M17").Value = 71.5
$ws.Range("N17").Value = -1000001228
$ws.Range("H129").Value = 1345.44
$ws.Range("I129").Value = 375.73334
$ws.Range("J129").Value = 2800
$ws.Range("K129").Value = 1127.20002
$ws.Range("L129").Value = 8400
$ws.Range("M129").Value = 3872.79998
$ws.Range("N129").Value = -18400
$ws.Range("H131").Value = 3771.5
$ws.Range("I131").Value = 1750
$ws.Range("J131").Value = 3996.111
$ws.Range("K131").Value = 5250
$ws.Range("L131").Value = 11988.333
$ws.Range("M131").Value = -210
$ws.Range("N131").Value = -22068.333

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H132").Value = 3814.4688
$ws.Range("I132").Value = 3492.25
$ws.Range("J132").Value = 6070
$ws.Range("K132").Value = 10476.75
$ws.Range("L132").Value = 18210
$ws.Range("M132").Value = -7946.75
$ws.Range("N132").Value = -23270
$ws.Range("H134").Value = 49997.332
$ws.Range("J134").Value = 49997.332
$ws.Range("L134").Value = 149991.996
$ws.Range("N134").Value = -155061.996

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 8815.192999999999
$ws.Range("J22").Value = 1485.7858
$ws.Range("L22").Value = 1485.7858
$ws.Range("N22").Value = -2075.7858
$ws.Range("H27").Value = 8815.192999999999
$ws.Range("J27").Value = 1485.7858
$ws.Range("L27").Value = 1485.7858
$ws.Range("N27").Value = -1699.7858
$ws.Range("H40").Value = 2963.4
$ws.Range("I40").Value = 2827.7646
$ws.Range("J40").Value = 3732
$ws.Range("K40").Value = 2827.7646
$ws.Range("L40").Value = 3732
$ws.Range("M40").Value = -2691.7646
$ws.Range("N40").Value = -4004
$ws.Range("H97").Value = 28167
$ws.Range("J97").Value = 28167
$ws.Range("L97").Value = 28167
$ws.Range("N97").Value = -30149
$ws.Range("H122").Value = 77976.57000000001
$ws.Range("I122").Value = 6352.1113
$ws.Range("J122").Value = 206900.6
$ws.Range("K122").Value = 19056.3339
$ws.Range("L122").Value = 620701.8
$ws.Range("M122").Value = -16606.3339
$ws.Range("N122").Value = -625601.8
$ws.Range("H132").Value = 58847.4
$ws.Range("I132").Value = 58847.4
$ws.Range("K132").Value = 176542.2
$ws.Range("M132").Value = -174012.2

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H81").Value = 108407
$ws.Range("J81").Value = 1999
$ws.Range("L81").Value = 3998
$ws.Range("N81").Value = -6120
$ws.Range("H84").Value = 108407
$ws.Range("J84").Value = 1999
$ws.Range("L84").Value = 19990
$ws.Range("N84").Value = -30598
$ws.Range("H107").Value = 4724.7144
$ws.Range("I107").Value = 3018.25
$ws.Range("J107").Value = 7000
$ws.Range("K107").Value = 9054.75
$ws.Range("L107").Value = 21000
$ws.Range("M107").Value = -7134.75
$ws.Range("N107").Value = -24840
$ws.Range("H132").Value = 3641.7917
$ws.Range("I132").Value = 2350.125
$ws.Range("J132").Value = 10100.125
$ws.Range("K132").Value = 7050.375
$ws.Range("L132").Value = 30300.375
$ws.Range("M132").Value = -4520.375
$ws.Range("N132").Value = -35360.375
